# CIDC-1278 first pass at redone docs
# Renames the "TCRseq Analysis" tab/strings to "TCR Analysis"/"TCR Runs",
# switches the active tab to "Legend", and tidies up the per-sheet
# selections that LibreOffice/Excel persist in sheetViews.

$wb = $excel.ActiveWorkbook

$wsAnalysis   = $wb.Worksheets.Item(1)   # "TCRseq Analysis" -> "TCR Analysis"
$wsExcluded   = $wb.Worksheets.Item(2)   # "Excluded Samples"
$wsLegend     = $wb.Worksheets.Item(3)   # "Legend"
$wsDataDict   = $wb.Worksheets.Item(4)   # "Data Dictionary"

# 1. Rename the first sheet (also drives the "TCRseq Analysis" -> "TCR Analysis"
#    text wherever the engine auto-links, but the shared-string copies in the
#    Legend/Data Dictionary text need to be edited explicitly below too).
$wsAnalysis.Name = "TCR Analysis"

# 2. Update the shared-string text that referenced the old tab/section names.
$wsAnalysis.Range("B7").Value2 = "TCR Runs"
$wsLegend.Range("B2").Value2 = "Legend for tab 'TCR Analysis'"
$wsLegend.Range("B7").Value2 = "Section 'TCR Runs' of tab 'TCR Analysis'"

# 3. Row-height tweaks on the Legend sheet that follow from the text edits
#    above (row 2 shrinks slightly, row 7 grows to fit the longer string).
$wsLegend.Rows.Item(2).RowHeight = 13.8
$wsLegend.Rows.Item(7).RowHeight = 23.95

# 4. Per-sheet selection bookkeeping.
$null = $wsAnalysis.Range("B8").Select()
$null = $wsExcluded.Range("B3").Select()
$null = $wsDataDict.Range("A1").Select()
$null = $wsLegend.Range("B8").Select()

# 5. Make "Legend" the active/visible tab (was "TCRseq Analysis").
$null = $wsLegend.Activate()

# 6. Best-effort: restore the "Normal 2"/"Normal 3" cell-style names to their
#    swapped order (cosmetic only - both styles resolve to the same format).
$styles = $wb.Styles
for ($i = 1; $i -le $styles.Count; $i++) {
    $style = $styles.Item($i)
    if ($style.Name -eq "Normal 3") {
        $style.Name = "__tmp_style_swap__"
    }
}
for ($i = 1; $i -le $styles.Count; $i++) {
    $style = $styles.Item($i)
    if ($style.Name -eq "Normal 2") {
        $style.Name = "Normal 3"
    }
}
for ($i = 1; $i -le $styles.Count; $i++) {
    $style = $styles.Item($i)
    if ($style.Name -eq "__tmp_style_swap__") {
        $style.Name = "Normal 2"
    }
}
